# Update cryptocurrency price/volume data to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '64.246.07'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = "'" + '3.491.87'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'" + '586.56'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = "'" + '134.04'
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D7').Value = "'" + '3.491.89'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'" + '0.485'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('D13').Value = "'" + '4.090.63'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = "'" + '0.119'
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = "'" + '3.494.69'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = "'" + '64.315.64'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = "'" + '25.17'
$ws.Range('E18').Value = '  -9.57%  '
$ws.Range('D19').Value = "'" + '9.85'
$ws.Range('E19').Value = '  -1.35%  '
$ws.Range('D20').Value = "'" + '5.73'
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('D21').Value = "'" + '13.59'
$ws.Range('E21').Value = '  -6.27%  '
$ws.Range('D22').Value = "'" + '387.47'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').Value = "'" + '3.634.48'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = "'" + '0.564'
$ws.Range('E24').Value = '  -2.49%  '
$ws.Range('D25').Value = "'" + '74.35'
$ws.Range('E25').Value = '  +1.98%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = "'" + '5.70'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = "'" + '1.00'
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = "'" + '7.37'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = "'" + '1.52'
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = "'" + '2.23'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = "'" + '8.25'
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').Value = "'" + '3.516.54'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('D36').Value = "'" + '0.148'
$ws.Range('E36').Value = '  +2.49%  '
$ws.Range('D37').Value = "'" + '23.46'
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('D38').Value = "'" + '5.23'
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').Value = "'" + '6.87'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').Value = "'" + '162.21'
$ws.Range('E41').Value = '  -2.72%  '
$ws.Range('D42').Value = "'" + '0.0781'
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('D43').Value = "'" + '0.804'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = "'" + '25.47'
$ws.Range('E45').Value = '  -6.27%  '
$ws.Range('D46').Value = "'" + '41.75'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('E48').Value = '  -2.40%  '
$ws.Range('D49').Value = "'" + '1.65'
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('D50').Value = "'" + '2.471.74'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').Value = "'" + '6.74'
$ws.Range('E51').Value = '  -2.26%  '
